$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 8 new rows before row 13 to make room for the new stretch-task rows ---
# (this shifts the former rows 13-40 down to rows 21-48, carrying their
#  formulas/styles/content with them, and Excel auto-adjusts the SUM/SUMIF
#  formulas that reference rows below the insertion point)
$ws.Rows.Item(13).Resize(8, 1).Insert()

# --- Fill in the new stretch task rows (13-20) ---

# Remove ingredients used to cook recipe
$ws.Range("A13").Value2 = "Remove ingredients used to cook recipe"
$ws.Range("B13").Value2 = "Complete functionality to remove ingredients used to cook recipe from pantry (web)"
$ws.Range("C13").Value2 = 0.5
$ws.Range("D13").Value2 = "Destiny"

$ws.Range("A14").Value2 = "Remove ingredients used to cook recipe"
$ws.Range("B14").Value2 = "Complete functionality to remove ingredients used to cook recipe from pantry (desktop)"
$ws.Range("C14").Value2 = 0.5
$ws.Range("D14").Value2 = "Janera"

# Add new recipe
$ws.Range("A15").Value2 = "Add new recipe"
$ws.Range("B15").Value2 = "Complete functionality add an new recipe (web)"
$ws.Range("C15").Value2 = 1
$ws.Range("D15").Value2 = "Matthew"

$ws.Range("A16").Value2 = "Add new recipe"
$ws.Range("B16").Value2 = "Complete functionality add an new recipe (desktop)"
$ws.Range("C16").Value2 = 0.5
$ws.Range("D16").Value2 = "Matthew"

# Share a recipe
$ws.Range("A17").Value2 = "Share a recipe"

$ws.Range("A18").Value2 = "Share a recipe"
$ws.Range("B18").Value2 = "Complete functionality to share a recipe (web)"
$ws.Range("C18").Value2 = 0.5
$ws.Range("D18").Value2 = "Destiny"

$ws.Range("B17").Value2 = "Complete functionality to share a recipe (desktop)"

$ws.Range("A19").Value2 = "Share a recipe"
$ws.Range("B19").Value2 = "Create UI to view shared recipes (desktop)"

$ws.Range("A20").Value2 = "Share a recipe"
$ws.Range("B20").Value2 = "Create UI to view shared recipes (web)"
$ws.Range("C20").Value2 = 1
$ws.Range("D20").Value2 = "Destiny"
$ws.Rows.Item(20).RowHeight = 14.25

# --- Demo Feedback row (reuses the row that used to be a blank "Feedback" row) ---
$ws.Range("A24").Value2 = "Demo Feedback"
$ws.Range("B24").Value2 = "Fix issues with desktop application mentioned in demo"
$ws.Range("C24").Value2 = 1
$ws.Range("D24").Value2 = "Janera"

# --- Update the Totals row (now row 28) to sum the full C3:C27 ranges ---
$ws.Range("C28").Formula = "=SUM(C3:C27)"
$ws.Range("E28").Formula = "=SUM(E3:E27)"
$ws.Range("H28").Formula = "=SUM(H3:H27)"
$ws.Range("I28").Formula = "=SUM(I3:I27)"
$ws.Rows.Item(28).RowHeight = 19.5

# --- Update the per-person SUMIF breakdown (now rows 33-35) ---
$ws.Range("C33").Formula = "=SUMIF(D3:D27, ""Matthew"", C3:C27)"
$ws.Range("D33").Formula = "=SUMIF(F3:F27, ""Matthew"", C3:C27)"
$ws.Range("C34").Formula = "=SUMIF(D3:D27, ""Destiny"", C3:C27)"
$ws.Range("D34").Formula = "=SUMIF(F3:F27, ""Destiny"", C3:C27)"
$ws.Range("C35").Formula = "=SUMIF(D3:D27, ""Janera"", C3:C27)"
$ws.Range("D35").Formula = "=SUMIF(F3:F27, ""Janera"", C3:C27)"

# --- Update the filter database defined name to the new totals-row block ---
$wb.Names("_xlnm._FilterDatabase").Delete()
$ws.Range("B28:I31").Name = "_xlnm._FilterDatabase"

# --- Update the current selection / view ---
$ws.Range("I28").Select()

# Nudge the worksheet's tracked used-range/dimension out to row 48 (the
# trailing formatting-only rows) without altering their look.
$ws.Range("A47").Font.Bold = $false
$ws.Range("A48").Font.Bold = $false
